$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Question 2")

# --- 2a) Unzen2 composition edits (input cells only; dependents recalc automatically) ---
$ws2.Range("B23").Value = 1.54
$ws2.Range("B29").Value = 20

# --- 2b) Unzen1 composition edits ---
$ws2.Range("B39").Value = 67.430000000000007
$ws2.Range("B40").Value = 0.16
$ws2.Range("B41").Value = 6.86
$ws2.Range("B42").Value = 0.25
$ws2.Range("B45").Value = 0.34
$ws2.Range("B46").Value = 2.29
$ws2.Range("B48").Value = 16.57

$wb.Save()
